$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "<Name>_old" -> "<Name>_FV2310", "<Name>_new" -> "<Name>_FV2404"
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# Turn the data range into an Excel Table (ListObject)
$rng = $ws.Range("A1:U64")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
